$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.800.98"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$ws.Range("D3").Value = "2.291.55"
$ws.Range("E3").Value = "  -1.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'299.52"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6
$ws.Range("D6").Value = "'96.44"
$ws.Range("E6").Value = "  -2.90%  "

# Row 7
$ws.Range("E7").Value = "  +0.30%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -3.10%  "

# Row 10
$ws.Range("D10").Value = "'35.66"
$ws.Range("E10").Value = "  -1.21%  "

# Row 11
$ws.Range("D11").Value = "'0.0786"
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("E12").Value = "  +0.83%  "

# Row 13
$ws.Range("D13").Value = "'17.63"
$ws.Range("E13").Value = "  -0.86%  "

# Row 14
$ws.Range("E14").Value = "  -2.23%  "

# Row 15
$ws.Range("D15").Value = "2.647.62"
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("D16").Value = "2.274.31"
$ws.Range("E16").Value = "  -1.09%  "

# Row 17
$ws.Range("D17").Value = "'0.773"
$ws.Range("E17").Value = "  -2.34%  "

# Row 18
$ws.Range("D18").Value = "42.749.85"
$ws.Range("E18").Value = "  -0.56%  "

# Row 19
$ws.Range("D19").Value = "'12.55"
$ws.Range("E19").Value = "  -5.22%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  -0.58%  "

# Row 21
$ws.Range("D21").Value = "'6.05"
$ws.Range("E21").Value = "  -2.30%  "

# Row 22
$ws.Range("D22").Value = "'67.76"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23
$ws.Range("D23").Value = "'242.15"
$ws.Range("E23").Value = "  +0.52%  "

# Row 24
$ws.Range("E24").Value = "  -1.37%  "

# Row 25
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -1.61%  "

# Row 27
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'4.01"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").Value = "'25.12"
$ws.Range("E28").Value = "  -0.13%  "

# Row 29
$ws.Range("D29").Value = "'165.89"
$ws.Range("E29").Value = "  -2.08%  "

# Row 30
$ws.Range("E30").Value = "  -0.75%  "

# Row 31
$ws.Range("E31").Value = "  -1.64%  "

# Row 32
$ws.Range("D32").Value = "'32.84"
$ws.Range("E32").Value = "  -1.67%  "

# Row 33
$ws.Range("E33").Value = "  +0.15%  "

# Row 34
$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "  -3.38%  "

# Row 35
$ws.Range("D35").Value = "'5.00"
$ws.Range("E35").Value = "  -3.80%  "

# Row 36
$ws.Range("D36").Value = "'17.08"
$ws.Range("E36").Value = "  -6.96%  "

# Row 37
$ws.Range("E37").Value = "  -1.31%  "

# Row 38
$ws.Range("D38").Value = "'0.0684"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39
$ws.Range("D39").Value = "'0.101"
$ws.Range("E39").Value = "  -1.73%  "

# Row 40
$ws.Range("E40").Value = "  -3.80%  "

# Row 41
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.74"
$ws.Range("E41").Value = "  -0.73%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.110"
$ws.Range("E42").Value = "  -0.06%  "

# Row 43
$ws.Range("D43").Value = "2.002.77"
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("E44").Value = "  -2.73%  "

# Row 45
$ws.Range("D45").Value = "'10.08"
$ws.Range("E45").Value = "  -0.65%  "

# Row 46
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
$ws.Range("D47").Value = "'17.04"
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("D48").Value = "'2.77"
$ws.Range("E48").Value = "  -2.51%  "

# Row 49
$ws.Range("D49").Value = "2.514.46"
$ws.Range("E49").Value = "  -1.12%  "

# Row 50
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.82"
$ws.Range("E50").Value = "  -4.68%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'53.11"
$ws.Range("E51").Value = "  -3.30%  "
